$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the status-check timestamp in the header cell F1
$ws.Range("F1").Value = "Last status check on: 17.01.2022 14:30"

# D10: was a text string "+0.2", now a real number
$ws.Range("D10").Value = 0.2

# E10: was a text string "2022-01-17 14:15:27", now a real datetime value
# formatted the same way as the other rows in column E
$ws.Range("E10").Value = 44578.5940625
$ws.Range("E10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
